# Auto-generated edit script: updates odds values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.96
$ws.Range("I5").Value = 4.33
$ws.Range("S5").Value = 1.62
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("X5").Value = 8
$ws.Range("AU5").Value = 9.5

# Row 7
$ws.Range("G7").Value = 1.39
$ws.Range("H7").Value = 4.55
$ws.Range("I7").Value = 6.6
$ws.Range("J7").Value = 1.85
$ws.Range("K7").Value = 2.42
$ws.Range("L7").Value = 6
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.82
$ws.Range("W7").Value = 7.8
$ws.Range("X7").Value = 7.1
$ws.Range("Y7").Value = 8.25
$ws.Range("Z7").Value = 9.25
$ws.Range("AA7").Value = 11
$ws.Range("AD7").Value = 9.25
$ws.Range("AE7").Value = 18.5
$ws.Range("AF7").Value = 75
$ws.Range("AG7").Value = 20
$ws.Range("AI7").Value = 21
$ws.Range("AL7").Value = 60
$ws.Range("AN7").Value = 3.3
$ws.Range("AO7").Value = 6.2
$ws.Range("AQ7").Value = 16.5
$ws.Range("AT7").Value = 3.2
$ws.Range("AU7").Value = 7.9
$ws.Range("AV7").Value = 65
$ws.Range("AW7").Value = 7.9
$ws.Range("AX7").Value = 35
$ws.Range("AY7").Value = 35

# Row 8
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 1.93
$ws.Range("L8").Value = 2.52
$ws.Range("W8").Value = 11.75
$ws.Range("X8").Value = 22
$ws.Range("Y8").Value = 12
$ws.Range("AB8").Value = 35
$ws.Range("AD8").Value = 6.5
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 55
$ws.Range("AJ8").Value = 17.5
$ws.Range("AL8").Value = 24
$ws.Range("AM8").Value = 400
$ws.Range("AP8").Value = 25
$ws.Range("AQ8").Value = 100
$ws.Range("AR8").Value = 120
$ws.Range("AS8").Value = 300
$ws.Range("AU8").Value = 6.8
$ws.Range("AW8").Value = 3.85
$ws.Range("AZ8").Value = 37
$ws.Range("BB8").Value = 200

# Row 13
$ws.Range("I13").Value = 2.88
$ws.Range("J13").Value = 3.2
$ws.Range("L13").Value = 3.6
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 8
$ws.Range("Q13").Value = 2.2
$ws.Range("R13").Value = 1.65
$ws.Range("S13").Value = 1.5
$ws.Range("T13").Value = 2.5
$ws.Range("AK13").Value = 26
$ws.Range("AT13").Value = 2.5
$ws.Range("AU13").Value = 8.5
$ws.Range("AY13").Value = 29

# Row 19
$ws.Range("G19").Value = 2.2
$ws.Range("I19").Value = 3.5
$ws.Range("J19").Value = 2.88
$ws.Range("L19").Value = 4
$ws.Range("U19").Value = 1.83
$ws.Range("V19").Value = 1.83
$ws.Range("X19").Value = 10
$ws.Range("AG19").Value = 10
$ws.Range("AI19").Value = 13
$ws.Range("AO19").Value = 12
$ws.Range("AZ19").Value = 67

# Row 20
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 3
$ws.Range("L20").Value = 4.33
$ws.Range("Q20").Value = 2.25
$ws.Range("R20").Value = 1.62
$ws.Range("X20").Value = 9
$ws.Range("AB20").Value = 29
$ws.Range("AC20").Value = 8
$ws.Range("AK20").Value = 34
$ws.Range("AM20").Value = 301
$ws.Range("AX20").Value = 21

# Row 26
$ws.Range("G26").Value = 2.9
$ws.Range("H26").Value = 2.67
$ws.Range("S26").Value = 1.53
$ws.Range("T26").Value = 2.2
$ws.Range("U26").Value = 1.91
$ws.Range("AD26").Value = 5.3
$ws.Range("AE26").Value = 15.5
$ws.Range("AF26").Value = 90
$ws.Range("AI26").Value = 10.25
$ws.Range("AJ26").Value = 35
$ws.Range("AK26").Value = 28
$ws.Range("AO26").Value = 16.5
$ws.Range("AT26").Value = 2.18

# Row 27
$ws.Range("G27").Value = 2.1
$ws.Range("I27").Value = 3.95
$ws.Range("J27").Value = 2.75
$ws.Range("U27").Value = 1.93
$ws.Range("W27").Value = 5.7
$ws.Range("X27").Value = 9
$ws.Range("AB27").Value = 37
$ws.Range("AC27").Value = 6.4
$ws.Range("AE27").Value = 15.5
$ws.Range("AG27").Value = 9
$ws.Range("AH27").Value = 21
$ws.Range("AN27").Value = 3.8
$ws.Range("AO27").Value = 11.25
$ws.Range("AP27").Value = 21
$ws.Range("AR27").Value = 90
$ws.Range("BB27").Value = 400

$wb.Save()